$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 9 days,
# and zero out all the notified production values in column B (rows 2-97).
for ($r = 2; $r -le 97; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value2 = $cellA.Value2 + 9

    $ws.Cells.Item($r, 2).Value2 = 0
}
